# edit.ps1 - apply the "React conversion" edit described by the diff
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Body: add a new paragraph run of text just before the trailing
#    bookmark (_GoBack) that sits in the last (empty) paragraph.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$insertRange = $lastPara.Range
$insertRange.Collapse(1)   # wdCollapseStart = 1
$insertRange.InsertBefore("In a React App we usually have a lot of JSX.  Browsers don" + [char]0x2019 + "t understand JSX, so we need to have a JavaScript compiler, such as " + [char]0x201C + "Babel" + [char]0x201D + " to transform our JSX code into regular JavaScript.")

# ---------------------------------------------------------------------
# 2. Header: collapse the "NextJS" + " Info - Page " runs (currently
#    split apart by proofErr spell-check markers) into a single run's
#    worth of text "NextJS Info - Page ".
# ---------------------------------------------------------------------
$hdr = $d.Sections.First.Headers.Item(1)
$hdr.Range.Find.Execute("NextJS Info - Page ", $false, $false, $false, $false, $false, $true, 1, $false, "NextJS###Info - Page ", 2)
$hdr.Range.Find.Execute("NextJS###Info - Page ", $false, $false, $false, $false, $false, $true, 1, $false, "NextJS Info - Page ", 2)
